$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I and J, matching style of existing headers (e.g. H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate data for columns I (I0) and J (IF) for rows 2-59
$data = New-Object 'object[,]' 58,2
$data[0,0] = 6; $data[0,1] = 7
$data[1,0] = 5; $data[1,1] = 6
$data[2,0] = 5; $data[2,1] = 6
$data[3,0] = 7; $data[3,1] = 7
$data[4,0] = 9; $data[4,1] = 9
$data[5,0] = 6; $data[5,1] = 6
$data[6,0] = 7; $data[6,1] = 7
$data[7,0] = 8; $data[7,1] = 8
$data[8,0] = 9; $data[8,1] = 9
$data[9,0] = 7; $data[9,1] = 7
$data[10,0] = 6; $data[10,1] = 6
$data[11,0] = 7; $data[11,1] = 7
$data[12,0] = 6; $data[12,1] = 6
$data[13,0] = 8; $data[13,1] = 8
$data[14,0] = 8; $data[14,1] = 8
$data[15,0] = 9; $data[15,1] = 9
$data[16,0] = 6; $data[16,1] = 6
$data[17,0] = 6; $data[17,1] = 6
$data[18,0] = 6; $data[18,1] = 6
$data[19,0] = 8; $data[19,1] = 8
$data[20,0] = 7; $data[20,1] = 7
$data[21,0] = 7; $data[21,1] = 8
$data[22,0] = 5; $data[22,1] = 6
$data[23,0] = 9; $data[23,1] = 9
$data[24,0] = 9; $data[24,1] = 9
$data[25,0] = 7; $data[25,1] = 7
$data[26,0] = 7; $data[26,1] = 7
$data[27,0] = 7; $data[27,1] = 8
$data[28,0] = 7; $data[28,1] = 7
$data[29,0] = 7; $data[29,1] = 7
$data[30,0] = 8; $data[30,1] = 8
$data[31,0] = 7; $data[31,1] = 7
$data[32,0] = 7; $data[32,1] = 7
$data[33,0] = 7; $data[33,1] = 7
$data[34,0] = 8; $data[34,1] = 8
$data[35,0] = 9; $data[35,1] = 9
$data[36,0] = 9; $data[36,1] = 9
$data[37,0] = 8; $data[37,1] = 9
$data[38,0] = 5; $data[38,1] = 6
$data[39,0] = 8; $data[39,1] = 9
$data[40,0] = 7; $data[40,1] = 7
$data[41,0] = 8; $data[41,1] = 8
$data[42,0] = 6; $data[42,1] = 6
$data[43,0] = 11; $data[43,1] = 11
$data[44,0] = 6; $data[44,1] = 6
$data[45,0] = 9; $data[45,1] = 10
$data[46,0] = 9; $data[46,1] = 9
$data[47,0] = 8; $data[47,1] = 8
$data[48,0] = 8; $data[48,1] = 8
$data[49,0] = 7; $data[49,1] = 7
$data[50,0] = 6; $data[50,1] = 6
$data[51,0] = 7; $data[51,1] = 7
$data[52,0] = 2; $data[52,1] = 2
$data[53,0] = 6; $data[53,1] = 6
$data[54,0] = 7; $data[54,1] = 7
$data[55,0] = 4; $data[55,1] = 5
$data[56,0] = 4; $data[56,1] = 4
$data[57,0] = 3; $data[57,1] = 3
$ws.Range("I2:J59").Value2 = $data
